# This workbook ("...Makerepayment1.xlsx") gets a new blank column inserted
# into the "Repayment Schedule" sheet (between the existing "In Advance" and
# "Late" columns), and the active sheet/selection bookmarks recorded in the
# workbook are updated to reflect where the user left off after editing
# (cursor ends up on the "Repayment Schedule" sheet, cell T10; the
# "Transactions" sheet keeps a plain cell selection instead of being the
# active tab with a full-sheet selection).

$wb = $excel.ActiveWorkbook

# --- Transactions sheet: drop the old full-sheet selection / active-tab
#     flag, replace with a simple selected cell. ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate() | Out-Null
$wsTransactions.Range("F24").Select() | Out-Null

# --- Repayment Schedule sheet: insert a new empty column before the old
#     column N ("Late"), shifting "Late" and "Over Due" one column to the
#     right (N->O, O->P, P->Q). The inserted column inherits formatting
#     from its neighbour (column M, "In Advance"). ---
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Activate() | Out-Null
$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# Leave this as the active sheet, with the cursor parked on T10, matching
# where editing finished.
$wsSchedule.Range("T10").Select() | Out-Null
